# Remove the bullet point about compile-time MSVC/GCC/Boost version checks.
# (commit message: "Boost version check.")
#
# The paragraph "Compile time checks to ensure MSVC, GCC, Boost, etc meet
# the minimum requirements." is deleted in its entirety, including its
# paragraph mark, so that the preceding and following bullets become
# adjacent list items.

$d = $word.ActiveDocument

$searchText = "Compile time checks to ensure MSVC, GCC, Boost, etc meet the minimum requirements."

$matchRange = $d.Content
$found = $matchRange.Find.Execute($searchText, $true, $false, $false, $false,
                                   $false, $true, 1, $false, "", 0)

if ($found) {
    # Build a fresh range spanning the matched sentence plus its trailing
    # paragraph mark, then delete it so the paragraph is removed (and the
    # surrounding bullets merge together) rather than left blank.
    $delRange = $d.Range($matchRange.Start, $matchRange.End + 1)
    $delRange.Delete()
}
